$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Professor Newman's notes..." paragraph: drop the proofErr markers
#    around "Matlab" by replacing the whole sentence with itself.
# ---------------------------------------------------------------------------
$old1 = "Professor Newman" + [char]8217 + "s notes on radial basis function networks describe a method for biasing the alpha nodes of a network of this type. However, his notes make the assumption that the inputs to the system are distributed across the range -1 to 1. Since these inputs do not have that distribution, I rescaled them using Equation 1 below. This rescaling could be incorporated into the weights of the alpha nodes, but for simplicity" + [char]8217 + "s sake, I just rescale the inputs at the beginning of my Matlab script."
$d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $old1, 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) "Once the inputs are rescaled..." paragraph: drop proofErr markers
#    around "Poressor" / "newman's".
# ---------------------------------------------------------------------------
$old2 = "Once the inputs are rescaled, the alpha node biases are generated as described in Poressor newman" + [char]8217 + "s notes. The weights from the inputs are chosen as to be random numbers evenly distributed from -1 to 1. Each bias weight is then chosen to be a random number evenly distributed"
$d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $old2, 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) "where i is an index..." sentence: drop proofErr markers around the
#    two "i" runs and "is".
# ---------------------------------------------------------------------------
$old3 = " where i is an index corresponding to an input (in this case i counts from 1 to 2) and w is the weight from that input into the alpha node. Figure 1 shows surface plot outputs of several alpha nodes, with scaled inputs superimposed in blue. "
$d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $old3, 2) | Out-Null

# ---------------------------------------------------------------------------
# 4) "Since the beta nodes have..." paragraph: drop proofErr markers around
#    "signum", "nalpha" (x2) and "seems".
# ---------------------------------------------------------------------------
$old4a = "For alpha nodes with a signum activation function, the beta node bias should be "
$d.Content.Find.Execute($old4a, $true, $false, $false, $false, $false, $true, 1, $false, $old4a, 2) | Out-Null

$old4b = "-n"
$d.Content.Find.Execute($old4b, $true, $false, $false, $false, $false, $true, 1, $false, $old4b, 2) | Out-Null

$old4c = "must be larger, and seems to be proportional to n"
$d.Content.Find.Execute($old4c, $true, $false, $false, $false, $false, $true, 1, $false, $old4c, 2) | Out-Null

$old4d = ". For 100 alpha nodes, 10 seems to be an appropriate number for "
$d.Content.Find.Execute($old4d, $true, $false, $false, $false, $false, $true, 1, $false, $old4d, 2) | Out-Null

# ---------------------------------------------------------------------------
# 5) "In order to smooth out..." paragraph: drop proofErr markers around
#    "softens".
# ---------------------------------------------------------------------------
$old5 = "In order to smooth out the beta node responses, I multiplied all of the input weights to the beta layer by a gain term. A gain less than one softens the response of the nodes. I found that a gain of about .05 resulted minimal error for a network with 25 beta nodes. Figure 3 shows that with these smoother beta responses, the output is also smoother and the output error goes down to about .022 (again, using optimal gamma weights)."
$d.Content.Find.Execute($old5, $true, $false, $false, $false, $false, $true, 1, $false, $old5, 2) | Out-Null

# ---------------------------------------------------------------------------
# 6) "Of course, increasing..." paragraph: drop proofErr markers around
#    "rms".
# ---------------------------------------------------------------------------
$old6 = "Of course, increasing the number of beta nodes also improves the quality of the fit, and reduces the need for the gain term in the beta weights. With 50 beta nodes, the best fit is achieved with a gain of .125, and the rms error drops to about .013. Figure 4 shows the responses of these beta nodes and the output response. Note the steep drop-off at the edges of the surface."
$d.Content.Find.Execute($old6, $true, $false, $false, $false, $false, $true, 1, $false, $old6, 2) | Out-Null

# ---------------------------------------------------------------------------
# 7) Add the new "Gamma Node Bias Selection by Random Perturbations" body
#    paragraphs, right after the heading and before the empty bookmark
#    paragraph.
# ---------------------------------------------------------------------------
$heading = $d.Content.Find.Execute("Gamma Node Bias Selection by Random Perturbations", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$headingPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "Gamma Node Bias Selection by Random Perturbations") {
        $headingPara = $p
        break
    }
}

$insertRange = $headingPara.Range
$insertRange.Collapse(0)

$delta = [char]916
$pm = [char]177
$rarr = [char]8594

$p1 = "As discussed above, the pseudoinverse method of tuning the gamma layer is not biologically plausible. An alternative (and biologically plausible) method of tuning is to initialize the gamma layer weights to some number (in this assignment 0) and then perform batch training. On each iteration, the weight vector Wgb is perturbed by adding to it a vector of random numbers, " + $delta + "Wgb. This " + $delta + "Wgb vector is generated using a uniform random distribution bounded by a constant " + $pm + " " + $delta + "max. The perturbed weights are used to simulate the entire set of inputs. If the rms error of the output decreases since the last iteration, the perturbed weights become the new weights. If the rms error increases, the perturbation is rejected. As the number of iterations approaches infinity, the rms error converges on the error of the optimal (pseudo inverse) solution, and all perturbations are rejected."

$p2 = "The parameter to be tuned in this training method is the bound of the random distribution " + $delta + "max. Choosing a small " + $delta + "max causes the network to approach the optimal solution slowly, but a larger " + $delta + "max will causes more of the perturbations to be rejected, ultimately resulting in slower training. The ideal " + $delta + "max is determined by the scale of the output relative to the outputs of the beta layer. In order to test different values of " + $delta + "max, the script iterates until one of two conditions is met. The first stop condition is that the error is within .01 of the optimal error. This stop condition signifies convergence on a valid solution. The second stop condition is that a perturbation has not been accepted in at least 100 iterations. This stop condition signifies a failure."

$insertRange.InsertParagraphAfter()
$insertRange.Collapse(0)
$insertRange.InsertAfter($p1)

$insertRange.Collapse(0)
$insertRange.InsertParagraphAfter()
$insertRange.Collapse(0)
$insertRange.InsertAfter($p2)

Write-Output "done"
